$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7791313529014587
$ws.Range("B1").Value = 1.462818026542664
$ws.Range("C1").Value = 5.597004413604736
$ws.Range("D1").Value = 3.15851092338562
$ws.Range("E1").Value = 1.491694331169128
